$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1): new columns I ("I0") and J ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing headers (bold font, thin box
# border, centered/top aligned) for the two new header cells.
$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data cells (rows 2-26): new columns I and J values ---
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 6
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 4
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 5
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 6
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 6
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 7
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 6
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 6
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 6
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 5
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 4
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 4
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 5
$ws.Range("I17").Value = 1
$ws.Range("J17").Value = 3
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = 5
$ws.Range("I19").Value = 9
$ws.Range("J19").Value = 9
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 9
$ws.Range("I21").Value = 8
$ws.Range("J21").Value = 9
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 7
$ws.Range("I23").Value = 5
$ws.Range("J23").Value = 5
$ws.Range("I24").Value = 4
$ws.Range("J24").Value = 6
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 3
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = 1
